# "Upload a revised excel file" -- add a second revised note below the
# existing one, then move the selection past it (matching the recorded
# workbook view from the authored commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 already holds "before revised"; B3 gets the new note as a string value
# (Excel stores it in the shared-string table, same as B2).
$ws.Range("B3").Value = "Revised!"

# After typing into B3 and pressing Enter, Excel's selection lands on B4.
$ws.Range("B4").Select()

# The workbook window was resized in the authored session; mirror that on
# the window object (best-effort -- matches the recorded view state).
$win = $wb.Windows.Item(1)
$win.Width = 19170
$win.Height = 8970
